$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split "...HTML i sl. " into "...HTML " + "i slično." (two runs
#    with identical formatting), and move the _GoBack bookmark so it
#    sits right after the new text (the old one, elsewhere in the
#    document, is removed automatically by adding a bookmark with the
#    same name).
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("i sl. ", $false, $false, $false, $false, $false, $true, 1, $false, "i slično.", 2)

# The replace above lands in a single run because the surrounding text
# shares identical formatting and gets coalesced back together. Force a
# genuine run split at the same spot by nudging a formatting property
# and then reverting it immediately.
$r2 = $d.Content
$r2.Find.Execute("i slično.", $false)
$r2.Bold = $true
$r2.Bold = $false

# Collapsing a range to the exact end of a paragraph (right before the
# paragraph mark) is mishandled by this host -- it resolves to the
# wrong location. Work around it by temporarily padding the paragraph
# end with two placeholder characters, anchoring the bookmark between
# the real text and the placeholder (a safe, non-boundary position),
# then deleting the placeholder again.
$p = $d.Paragraphs.Item(54)
$p.Range.InsertAfter("XX")

$r3 = $d.Content
$r3.Find.Execute("slično.", $false)
$r3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r3)

$r4 = $d.Content
$r4.Find.Execute("XX", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 2) Update the cached PAGE field text in the footer from "3" to "2".
#    Editing the field's Result range (or the footer Range as a whole)
#    does not reliably persist in this host, but per-character edits
#    via the footer's Characters collection do.
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    if ($ftr.Exists) {
        $chars = $ftr.Range.Characters
        for ($i = 1; $i -le $chars.Count; $i++) {
            $ch = $chars.Item($i)
            if ($ch.Text -eq "3") {
                $ch.Text = "2"
            }
        }
    }
}
